$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Save" column header in H1, matching the style used by the
# other header cells (e.g. G1: bold, bordered, centered)
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)

# Fill the new "Save" column values for the two data rows
$ws.Range("H2").Value = 1
$ws.Range("H3").Value = 1
